# Auto-generated edit script.
# Applies the per-cell numeric updates described in the commit diff
# for Sheets/Asura_Profits.xlsx (market-price refresh across the 8 crafting-job sheets).
#
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3922.5386
$ws.Range("I74").Value = 3164.3333
$ws.Range("K74").Value = 3164.3333
$ws.Range("M74").Value = -2228.3333
$ws.Range("H76").Value = 3951.5557
$ws.Range("I76").Value = 3927.1538
$ws.Range("J76").Value = 4015
$ws.Range("K76").Value = 3927.1538
$ws.Range("L76").Value = 4015
$ws.Range("M76").Value = -3612.1538
$ws.Range("N76").Value = -4645
$ws.Range("H77").Value = 3922.5386
$ws.Range("I77").Value = 3164.3333
$ws.Range("K77").Value = 15821.6665
$ws.Range("M77").Value = -11141.6665
$ws.Range("H79").Value = 3951.5557
$ws.Range("I79").Value = 3927.1538
$ws.Range("J79").Value = 4015
$ws.Range("K79").Value = 3927.1538
$ws.Range("L79").Value = 4015
$ws.Range("M79").Value = -2835.1538
$ws.Range("N79").Value = -6199
$ws.Range("H112").Value = 2156.1333
$ws.Range("J112").Value = 2364.7693
$ws.Range("L112").Value = 7094.3079
$ws.Range("N112").Value = -9310.3079
$ws.Range("H116").Value = 25002736
$ws.Range("I116").Value = 33335764
$ws.Range("J116").Value = 3650
$ws.Range("K116").Value = 33335764
$ws.Range("L116").Value = 3650
$ws.Range("M116").Value = -33332322
$ws.Range("N116").Value = -10534
$ws.Range("H127").Value = 956.34
$ws.Range("I127").Value = 477.25
$ws.Range("K127").Value = 1431.75
$ws.Range("M127").Value = 3528.25

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3404.8572
$ws.Range("I63").Value = 3692.4614
$ws.Range("J63").Value = 2937.5
$ws.Range("K63").Value = 3692.4614
$ws.Range("L63").Value = 2937.5
$ws.Range("M63").Value = -3006.4614
$ws.Range("N63").Value = -4309.5
$ws.Range("H66").Value = 3404.8572
$ws.Range("I66").Value = 3692.4614
$ws.Range("J66").Value = 2937.5
$ws.Range("K66").Value = 18462.307
$ws.Range("L66").Value = 14687.5
$ws.Range("M66").Value = -15030.307
$ws.Range("N66").Value = -21551.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4166.9
$ws.Range("I105").Value = 3961.125
$ws.Range("K105").Value = 3961.125
$ws.Range("M105").Value = -2214.125

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9624.75
$ws.Range("J4").Value = 9624.75
$ws.Range("L4").Value = 9624.75
$ws.Range("N4").Value = -9848.75
$ws.Range("H62").Value = 86015
$ws.Range("I62").Value = 86015
$ws.Range("K62").Value = 86015
$ws.Range("M62").Value = -85391
$ws.Range("H65").Value = 86015
$ws.Range("I65").Value = 86015
$ws.Range("K65").Value = 430075
$ws.Range("M65").Value = -426955

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1259.6
$ws.Range("I4").Value = 99.333336
$ws.Range("J4").Value = 3000
$ws.Range("K4").Value = 298.000008
$ws.Range("L4").Value = 9000
$ws.Range("M4").Value = -186.000008
$ws.Range("N4").Value = -9224
$ws.Range("H22").Value = 20834604
$ws.Range("J22").Value = 1610
$ws.Range("L22").Value = 4830
$ws.Range("N22").Value = -5168
$ws.Range("H27").Value = 20834604
$ws.Range("J27").Value = 1610
$ws.Range("L27").Value = 4830
$ws.Range("N27").Value = -5034
$ws.Range("H131").Value = 3267.4036
$ws.Range("I131").Value = 14681.429
$ws.Range("J131").Value = 1669.44
$ws.Range("K131").Value = 44044.287
$ws.Range("L131").Value = 5008.32
$ws.Range("M131").Value = -39004.287
$ws.Range("N131").Value = -15088.32

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 4497.5
$ws.Range("I5").Value = 1996.6666
$ws.Range("K5").Value = 1996.6666
$ws.Range("M5").Value = -1884.6666
$ws.Range("H52").Value = 50000
$ws.Range("J52").Value = 50000
$ws.Range("L52").Value = 50000
$ws.Range("N52").Value = -50518
$ws.Range("H70").Value = 5884.3257
$ws.Range("I70").Value = 6183.273
$ws.Range("J70").Value = 5781.5625
$ws.Range("K70").Value = 6183.273
$ws.Range("L70").Value = 5781.5625
$ws.Range("M70").Value = -5913.273
$ws.Range("N70").Value = -6321.5625
$ws.Range("H73").Value = 5884.3257
$ws.Range("I73").Value = 6183.273
$ws.Range("J73").Value = 5781.5625
$ws.Range("K73").Value = 6183.273
$ws.Range("L73").Value = 5781.5625
$ws.Range("M73").Value = -5247.273
$ws.Range("N73").Value = -7653.5625
$ws.Range("H80").Value = 2791.5833
$ws.Range("I80").Value = 2714.2856
$ws.Range("J80").Value = 2899.8
$ws.Range("K80").Value = 2714.2856
$ws.Range("L80").Value = 2899.8
$ws.Range("M80").Value = -1716.2856
$ws.Range("N80").Value = -4895.8
$ws.Range("H83").Value = 2791.5833
$ws.Range("I83").Value = 2714.2856
$ws.Range("J83").Value = 2899.8
$ws.Range("K83").Value = 13571.428
$ws.Range("L83").Value = 14499
$ws.Range("M83").Value = -8579.428
$ws.Range("N83").Value = -24483

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9930.315000000001
$ws.Range("J2").Value = 9930.315000000001
$ws.Range("L2").Value = 9930.315000000001
$ws.Range("N2").Value = -10154.315
$ws.Range("H7").Value = 4427.5
$ws.Range("I7").Value = 5066.6665
$ws.Range("J7").Value = 4253.1816
$ws.Range("K7").Value = 5066.6665
$ws.Range("L7").Value = 4253.1816
$ws.Range("M7").Value = -4954.6665
$ws.Range("N7").Value = -4477.1816
$ws.Range("H122").Value = 13164105
$ws.Range("I122").Value = 22732364
$ws.Range("J122").Value = 7750
$ws.Range("K122").Value = 68197092
$ws.Range("L122").Value = 23250
$ws.Range("M122").Value = -68194642
$ws.Range("N122").Value = -28150
$ws.Range("H126").Value = 4427.5
$ws.Range("I126").Value = 5066.6665
$ws.Range("J126").Value = 4253.1816
$ws.Range("K126").Value = 15199.9995
$ws.Range("L126").Value = 12759.5448
$ws.Range("M126").Value = -12729.9995
$ws.Range("N126").Value = -17699.5448
$ws.Range("H132").Value = 5261.0713
$ws.Range("I132").Value = 5238.4287
$ws.Range("K132").Value = 15715.2861
$ws.Range("M132").Value = -13185.2861

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = $null
$ws.Range("H5").Value = 400001.34
$ws.Range("I5").Value = 1000000
$ws.Range("K5").Value = 1000000
$ws.Range("M5").Value = -999888

